# Remove the first 4 data rows (years 1985-1988, originally sheet rows 2-5)
# so that the remaining data shifts up and the sheet now spans A1:E38
# instead of A1:E42.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:E5").EntireRow.Delete()
